$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.467421333333334
$ws.Range("H2").Value = 25.402264
$ws.Range("I2").Value = 0.2732469334691616
$ws.Range("J2").Value = 0.312800300005396
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.68421466666667
$ws.Range("N2").Value = 32.052644
$ws.Range("O2").Value = 0.04705285980693976
$ws.Range("P2").Value = 0.04892736897547583
$ws.Range("Q2").Value = 90.46774719844623
$ws.Range("R2").Value = 814.2097247860161
$ws.Range("S2").Value = 0.01285704965320066
$ws.Range("T2").Value = 0.01530449569400354

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.467421333333334
$ws.Range("H3").Value = 25.402264
$ws.Range("I3").Value = 0.2732469334691616
$ws.Range("J3").Value = 0.312800300005396
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 70.36235166666667
$ws.Range("N3").Value = 211.087055
$ws.Range("O3").Value = 0.3098730203341347
$ws.Range("P3").Value = 0.3222178559101571
$ws.Range("Q3").Value = 595.7876775658357
$ws.Range("R3").Value = 5362.089098092521
$ws.Range("S3").Value = 0.08467185257112946
$ws.Range("T3").Value = 0.1007898419957926

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.467421333333334
$ws.Range("H4").Value = 25.402264
$ws.Range("I4").Value = 0.2732469334691616
$ws.Range("J4").Value = 0.312800300005396
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 59.09107466666666
$ws.Range("N4").Value = 177.273224
$ws.Range("O4").Value = 0.2602347611759026
$ws.Range("P4").Value = 0.2706020894912812
$ws.Range("Q4").Value = 500.3490262421262
$ws.Range("R4").Value = 4503.141236179136
$ws.Range("S4").Value = 0.071108350473395
$ws.Range("T4").Value = 0.08464441477495978

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.467421333333334
$ws.Range("H5").Value = 25.402264
$ws.Range("I5").Value = 0.2732469334691616
$ws.Range("J5").Value = 0.312800300005396
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 60.83231733333333
$ws.Range("N5").Value = 182.496952
$ws.Range("O5").Value = 0.2679031251727568
$ws.Range("P5").Value = 0.2785759485989269
$ws.Range("Q5").Value = 515.0928615443698
$ws.Range("R5").Value = 4635.835753899329
$ws.Range("S5").Value = 0.07320370742026075
$ws.Range("T5").Value = 0.0871386402960321

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.467421333333334
$ws.Range("H6").Value = 25.402264
$ws.Range("I6").Value = 0.2732469334691616
$ws.Range("J6").Value = 0.312800300005396
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 26.0983795
$ws.Range("N6").Value = 52.196759
$ws.Range("O6").Value = 0.1149362335102661
$ws.Range("P6").Value = 0.07967673702415903
$ws.Range("Q6").Value = 220.9859753437294
$ws.Range("R6").Value = 1325.915852062376
$ws.Range("S6").Value = 0.03140597335117571
$ws.Range("T6").Value = 0.02492290724460799

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3343473333333333
$ws.Range("H7").Value = 1.003042
$ws.Range("I7").Value = 0.01078951666043526
$ws.Range("J7").Value = 0.01235133366529898
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.68421466666667
$ws.Range("N7").Value = 32.052644
$ws.Range("O7").Value = 0.04705285980693976
$ws.Range("P7").Value = 0.04892736897547583
$ws.Range("Q7").Value = 3.572238682560889
$ws.Range("R7").Value = 32.150148143048
$ws.Range("S7").Value = 0.0005076776148081009
$ws.Range("T7").Value = 0.0006043182595812996

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.3343473333333333
$ws.Range("H8").Value = 1.003042
$ws.Range("I8").Value = 0.01078951666043526
$ws.Range("J8").Value = 0.01235133366529898
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 70.36235166666667
$ws.Range("N8").Value = 211.087055
$ws.Range("O8").Value = 0.3098730203341347
$ws.Range("P8").Value = 0.3222178559101571
$ws.Range("Q8").Value = 23.52546464681222
$ws.Range("R8").Value = 211.72918182131
$ws.Range("S8").Value = 0.003343380115514539
$ws.Range("T8").Value = 0.00397982025126358

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.3343473333333333
$ws.Range("H9").Value = 1.003042
$ws.Range("I9").Value = 0.01078951666043526
$ws.Range("J9").Value = 0.01235133366529898
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 59.09107466666666
$ws.Range("N9").Value = 177.273224
$ws.Range("O9").Value = 0.2602347611759026
$ws.Range("P9").Value = 0.2706020894912812
$ws.Range("Q9").Value = 19.75694323860089
$ws.Range("R9").Value = 177.812489147408
$ws.Range("S9").Value = 0.002807807291331791
$ws.Range("T9").Value = 0.00334229669783391

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3343473333333333
$ws.Range("H10").Value = 1.003042
$ws.Range("I10").Value = 0.01078951666043526
$ws.Range("J10").Value = 0.01235133366529898
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 60.83231733333333
$ws.Range("N10").Value = 182.496952
$ws.Range("O10").Value = 0.2679031251727568
$ws.Range("P10").Value = 0.2785759485989269
$ws.Range("Q10").Value = 20.33912308088711
$ws.Range("R10").Value = 183.052107727984
$ws.Range("S10").Value = 0.002890545232434131
$ws.Range("T10").Value = 0.003440784492272524

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3343473333333333
$ws.Range("H11").Value = 1.003042
$ws.Range("I11").Value = 0.01078951666043526
$ws.Range("J11").Value = 0.01235133366529898
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 26.0983795
$ws.Range("N11").Value = 52.196759
$ws.Range("O11").Value = 0.1149362335102661
$ws.Range("P11").Value = 0.07967673702415903
$ws.Range("Q11").Value = 8.725923590146333
$ws.Range("R11").Value = 52.355541540878
$ws.Range("S11").Value = 0.001240106406346694
$ws.Range("T11").Value = 0.0009841139643476692

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.140032333333333
$ws.Range("H12").Value = 12.420097
$ws.Range("I12").Value = 0.1336004309946363
$ws.Range("J12").Value = 0.1529395201819853
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.68421466666667
$ws.Range("N12").Value = 32.052644
$ws.Range("O12").Value = 0.04705285980693976
$ws.Range("P12").Value = 0.04892736897547583
$ws.Range("Q12").Value = 44.23299417627423
$ws.Range("R12").Value = 398.096947586468
$ws.Range("S12").Value = 0.006286282349737349
$ws.Range("T12").Value = 0.007482928334876226

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.140032333333333
$ws.Range("H13").Value = 12.420097
$ws.Range("I13").Value = 0.1336004309946363
$ws.Range("J13").Value = 0.1529395201819853
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 70.36235166666667
$ws.Range("N13").Value = 211.087055
$ws.Range("O13").Value = 0.3098730203341347
$ws.Range("P13").Value = 0.3222178559101571
$ws.Range("Q13").Value = 291.3024109493705
$ws.Range("R13").Value = 2621.721698544335
$ws.Range("S13").Value = 0.04139916907025008
$ws.Range("T13").Value = 0.04927984427696749

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.140032333333333
$ws.Range("H14").Value = 12.420097
$ws.Range("I14").Value = 0.1336004309946363
$ws.Range("J14").Value = 0.1529395201819853
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 59.09107466666666
$ws.Range("N14").Value = 177.273224
$ws.Range("O14").Value = 0.2602347611759026
$ws.Range("P14").Value = 0.2706020894912812
$ws.Range("Q14").Value = 244.6389597314142
$ws.Range("R14").Value = 2201.750637582728
$ws.Range("S14").Value = 0.03476747625288682
$ws.Range("T14").Value = 0.04138575372703919

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.140032333333333
$ws.Range("H15").Value = 12.420097
$ws.Range("I15").Value = 0.1336004309946363
$ws.Range("J15").Value = 0.1529395201819853
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 60.83231733333333
$ws.Range("N15").Value = 182.496952
$ws.Range("O15").Value = 0.2679031251727568
$ws.Range("P15").Value = 0.2785759485989269
$ws.Range("Q15").Value = 251.8477606715937
$ws.Range("R15").Value = 2266.629846044344
$ws.Range("S15").Value = 0.0357919729878903
$ws.Range("T15").Value = 0.04260527191296127

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.140032333333333
$ws.Range("H16").Value = 12.420097
$ws.Range("I16").Value = 0.1336004309946363
$ws.Range("J16").Value = 0.1529395201819853
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 26.0983795
$ws.Range("N16").Value = 52.196759
$ws.Range("O16").Value = 0.1149362335102661
$ws.Range("P16").Value = 0.07967673702415903
$ws.Range("Q16").Value = 108.0481349776038
$ws.Range("R16").Value = 648.288809865623
$ws.Range("S16").Value = 0.01535553033387171
$ws.Range("T16").Value = 0.0121857219301411

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.291073
$ws.Range("H17").Value = 18.873219
$ws.Range("I17").Value = 0.2030153381778063
$ws.Range("J17").Value = 0.2324024569332694
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 10.68421466666667
$ws.Range("N17").Value = 32.052644
$ws.Range("O17").Value = 0.04705285980693976
$ws.Range("P17").Value = 0.04892736897547583
$ws.Range("Q17").Value = 67.21517441567067
$ws.Range("R17").Value = 604.936569741036
$ws.Range("S17").Value = 0.009552452245938786
$ws.Range("T17").Value = 0.0113708407611812

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 6.291073
$ws.Range("H18").Value = 18.873219
$ws.Range("I18").Value = 0.2030153381778063
$ws.Range("J18").Value = 0.2324024569332694
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 70.36235166666667
$ws.Range("N18").Value = 211.087055
$ws.Range("O18").Value = 0.3098730203341347
$ws.Range("P18").Value = 0.3222178559101571
$ws.Range("Q18").Value = 442.6546907866717
$ws.Range("R18").Value = 3983.892217080045
$ws.Range("S18").Value = 0.06290897601531262
$ws.Range("T18").Value = 0.07488422138129067

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 6.291073
$ws.Range("H19").Value = 18.873219
$ws.Range("I19").Value = 0.2030153381778063
$ws.Range("J19").Value = 0.2324024569332694
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 59.09107466666666
$ws.Range("N19").Value = 177.273224
$ws.Range("O19").Value = 0.2602347611759026
$ws.Range("P19").Value = 0.2706020894912812
$ws.Range("Q19").Value = 371.7462643764507
$ws.Range("R19").Value = 3345.716379388056
$ws.Range("S19").Value = 0.05283164804574653
$ws.Range("T19").Value = 0.06288859044905018

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 6.291073
$ws.Range("H20").Value = 18.873219
$ws.Range("I20").Value = 0.2030153381778063
$ws.Range("J20").Value = 0.2324024569332694
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 60.83231733333333
$ws.Range("N20").Value = 182.496952
$ws.Range("O20").Value = 0.2679031251727568
$ws.Range("P20").Value = 0.2785759485989269
$ws.Range("Q20").Value = 382.7005491031653
$ws.Range("R20").Value = 3444.304941928488
$ws.Range("S20").Value = 0.0543884435558384
$ws.Range("T20").Value = 0.06474173489690677

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 6.291073
$ws.Range("H21").Value = 18.873219
$ws.Range("I21").Value = 0.2030153381778063
$ws.Range("J21").Value = 0.2324024569332694
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 26.0983795
$ws.Range("N21").Value = 52.196759
$ws.Range("O21").Value = 0.1149362335102661
$ws.Range("P21").Value = 0.07967673702415903
$ws.Range("Q21").Value = 164.1868106162035
$ws.Range("R21").Value = 985.1208636972209
$ws.Range("S21").Value = 0.02333381831497
$ws.Range("T21").Value = 0.01851706944484055

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 11.755292
$ws.Range("H22").Value = 23.510584
$ws.Range("I22").Value = 0.3793477806979606
$ws.Range("J22").Value = 0.2895063892140504
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 10.68421466666667
$ws.Range("N22").Value = 32.052644
$ws.Range("O22").Value = 0.04705285980693976
$ws.Range("P22").Value = 0.04892736897547583
$ws.Range("Q22").Value = 125.5960631973493
$ws.Range("R22").Value = 753.5763791840959
$ws.Range("S22").Value = 0.01784939794325487
$ws.Range("T22").Value = 0.01416478592583356

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 11.755292
$ws.Range("H23").Value = 23.510584
$ws.Range("I23").Value = 0.3793477806979606
$ws.Range("J23").Value = 0.2895063892140504
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 70.36235166666667
$ws.Range("N23").Value = 211.087055
$ws.Range("O23").Value = 0.3098730203341347
$ws.Range("P23").Value = 0.3222178559101571
$ws.Range("Q23").Value = 827.1299896483533
$ws.Range("R23").Value = 4962.779937890119
$ws.Range("S23").Value = 0.117549642561928
$ws.Range("T23").Value = 0.09328412800484275

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 11.755292
$ws.Range("H24").Value = 23.510584
$ws.Range("I24").Value = 0.3793477806979606
$ws.Range("J24").Value = 0.2895063892140504
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 59.09107466666666
$ws.Range("N24").Value = 177.273224
$ws.Range("O24").Value = 0.2602347611759026
$ws.Range("P24").Value = 0.2706020894912812
$ws.Range("Q24").Value = 694.6328373004692
$ws.Range("R24").Value = 4167.797023802816
$ws.Range("S24").Value = 0.09871947911254246
$ws.Range("T24").Value = 0.07834103384239817

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 11.755292
$ws.Range("H25").Value = 23.510584
$ws.Range("I25").Value = 0.3793477806979606
$ws.Range("J25").Value = 0.2895063892140504
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 60.83231733333333
$ws.Range("N25").Value = 182.496952
$ws.Range("O25").Value = 0.2679031251727568
$ws.Range("P25").Value = 0.2785759485989269
$ws.Range("Q25").Value = 715.1016532899946
$ws.Range("R25").Value = 4290.609919739967
$ws.Range("S25").Value = 0.1016284559763332
$ws.Range("T25").Value = 0.08064951700075422

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 11.755292
$ws.Range("H26").Value = 23.510584
$ws.Range("I26").Value = 0.3793477806979606
$ws.Range("J26").Value = 0.2895063892140504
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 26.0983795
$ws.Range("N26").Value = 52.196759
$ws.Range("O26").Value = 0.1149362335102661
$ws.Range("P26").Value = 0.07967673702415903
$ws.Range("Q26").Value = 306.794071749314
$ws.Range("R26").Value = 1227.176286997256
$ws.Range("S26").Value = 0.04360080510390203
$ws.Range("T26").Value = 0.02306692444022173
